# Move the SmartArt diagram graphic frame further down the slide.
# (Original <a:off x="1137600" y="1407600"/> -> <a:off x="1137600" y="2847600"/>,
#  i.e. the top offset grows from 1407600 EMU to 2847600 EMU.)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)

# PowerPoint's Shape.Top/.Left are expressed in points (1 pt = 12700 EMU).
$targetTopEmu = 2847600
$shp.Top = $targetTopEmu / 914400 * 72
